$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.914.68"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.83%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.296.26"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.24%  "

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.17%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "316.62"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.12%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "104.14"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.03%  "

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -1.00%  "

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.25%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.602"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -1.62%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.46"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -1.72%  "

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.72%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.47"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +1.43%  "

$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +2.31%  "

$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +4.42%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.35"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.24%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.644.06"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.22%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.298.14"

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "42.836.42"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +1.00%  "

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.03%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.85"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +26.32%  "

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.52%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "73.88"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.74%  "

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.65%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "263.30"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -4.71%  "

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -3.45%  "

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.13%  "

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.72%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.16"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +21.64%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.35"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.06%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "22.33"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -1.93%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "37.60"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +5.10%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "166.69"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +1.25%  "

$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -3.83%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.59"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.46%  "

$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -1.63%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.56"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.65%  "

$ws.Range("B38").NumberFormat = "@"
$ws.Range("B38").Value = "NEARProtocol"
$ws.Range("C38").NumberFormat = "@"
$ws.Range("C38").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.87"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +3.54%  "

$ws.Range("B39").NumberFormat = "@"
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").NumberFormat = "@"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0351"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -4.76%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.67"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -2.83%  "

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +5.30%  "

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +1.12%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "69.63"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.54%  "

$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.34%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "92.35"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.37%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.21"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +1.48%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "114.42"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +1.36%  "

$ws.Range("B48").NumberFormat = "@"
$ws.Range("B48").Value = "Maker"
$ws.Range("C48").NumberFormat = "@"
$ws.Range("C48").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.721.60"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +7.86%  "

$ws.Range("B49").NumberFormat = "@"
$ws.Range("B49").Value = "ordi"
$ws.Range("C49").NumberFormat = "@"
$ws.Range("C49").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "80.08"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -3.21%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.79"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -1.33%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "5.18"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +1.85%  "
